# Regen save_data to use K (strikeouts) instead of Strike# placeholder values.
# Update column G ("K") for the affected rows with the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 2
    24 = 2
    25 = 1
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
